# Applies the "added classifications of the sample runs" edit to
# Warehouse/testing/FF/judgement_logs_run4o.xlsx
#
# Summary of the change:
#   - Row 2 (file run4o_discovery_04_29_2025 at_07;03;59B.json) is replaced by the
#     run4o_discovery_06_27_2025 at_21;54;02N.json judgement -> "no_decision"
#   - Row 3's judgement is updated to the "no_decision" wording (filename unchanged)
#   - Row 4 is replaced by the run4o_discovery_04_29_2025 at_06;58;37N.json judgement -> "no_decision"
#   - A brand new row 5 is inserted holding the original "Barbie" / row4 content
#     (run4o_discovery_04_29_2025 at_07;00;29B.json -> "Barbie_was_selected")
#   - Four new rows (6-9) are appended with additional sample judgements
#   - The old row 5 (run4o_discovery_06_27_2025 at_22;14;40C.json, "both_movies") moves to row 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, extend the A-column "index" formatting (border/bold/centered style)
# down into the new rows 6:9 by copying the format already used by A2:A5. ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("B2").Value = "./Warehouse/testing/FF/run4o_discovery_06_27_2025 at_21;54;02N.json"
$ws.Range("C2").Value = "MSG: None`n`nMSG: I have recorded the decision reflecting that no consensus was reached regarding the movie to show on Friday.`n"
$ws.Range("D2").Value = "no_decision, "

# --- Row 3 (filename A3/B3 stays the same) ---
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that a consensus has not been reached on which movie to show on Friday.`n"
$ws.Range("D3").Value = "no_decision, "

# --- Row 4 ---
$ws.Range("B4").Value = "./Warehouse/testing/FF/run4o_discovery_04_29_2025 at_06;58;37N.json"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has been recorded as no decision being made.`n"
$ws.Range("D4").Value = "no_decision, "

# --- Row 5 (new) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "./Warehouse/testing/FF/run4o_discovery_04_29_2025 at_07;00;29B.json"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("D5").Value = "Barbie_was_selected, "

# --- Row 6 (new) ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "./Warehouse/testing/FF/run4o_discovery_05_04_2025 at_21;18;06N.json"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The conversation concluded without a decision on what movie to show on Friday, so no further action is required.`n"
$ws.Range("D6").Value = "no_decision, "

# --- Row 7 (new) ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "./Warehouse/testing/FF/run4o_discovery_04_29_2025 at_07;03;59B.json"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to not show a movie this Friday, as there was no consensus among the committee members.`n"
$ws.Range("D7").Value = "no_decision, "

# --- Row 8 (previously row 5) ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "./Warehouse/testing/FF/run4o_discovery_06_27_2025 at_22;14;40C.json"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made, so no action will be taken regarding acquiring movie rights.`n"
$ws.Range("D8").Value = "no_decision, "

# --- Row 9 (new) ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "./Warehouse/testing/FF/run4o_discovery_04_29_2025 at_07;05;46N.json"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("D9").Value = "both_movies, "
